# Applies the "Updated all sumsum.txt and excel files" commit to the workbook.
# Two worksheets: "table&first graph" (sheet1) and "bounds in mem or regs graph" (sheet2).
# Only the raw input cells are changed below - every other cell in the workbook is a
# formula (direct reference, shared formula or OFFSET-based) that recalculates from
# these inputs automatically.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("table&first graph")
$ws2 = $wb.Worksheets.Item("bounds in mem or regs graph")

# ---------------------------------------------------------------------------
# Sheet 1 ("table&first graph") raw data edits
# ---------------------------------------------------------------------------

# Header: benchmark name "coremk" -> "coremk_ch2"
$ws1.Range("K2").Value = "coremk_ch2"

# Row 4 (unsafe, 9 register pairs)
$ws1.Range("K4").Value = 58.9
$ws1.Range("N4").Value = 70.2
$ws1.Range("O4").Value = 67

# Row 5 (unsafe, 11 register pairs)
$ws1.Range("K5").Value = 46.7
$ws1.Range("N5").Value = 95.4
$ws1.Range("O5").Value = 77.4

# Row 7 (safe writes, 11 register pairs)
$ws1.Range("K7").Value = 76.7
$ws1.Range("N7").Value = 88.2
$ws1.Range("O7").Value = 104.6

# Row 8 (safe reads and writes, 11 register pairs)
$ws1.Range("K8").Value = 50.5
$ws1.Range("N8").Value = 102.6
$ws1.Range("O8").Value = 81.6

# Row 10 (safe writes, bounds in mem)
$ws1.Range("K10").Value = 155
$ws1.Range("N10").Value = 120.8
$ws1.Range("O10").Value = 173.9

# Row 11 (safe reads and writes, bounds in mem)
$ws1.Range("K11").Value = 58.2
$ws1.Range("N11").Value = 106.2
$ws1.Range("O11").Value = 88.1

# Row 13 (03) Obj/arr st)
$ws1.Range("K13").Value = 11.7
$ws1.Range("N13").Value = 8.8
$ws1.Range("O13").Value = 10.3

# Note text at Q5:R5 - content and highlight colour (orange -> red)
$ws1.Range("Q5").Value = "UPDATED 20180326"
$ws1.Range("Q5:R5").Interior.Color = RGB(255, 0, 0)
$ws1.Range("R5").Font.Color = RGB(255, 0, 0)

# ---------------------------------------------------------------------------
# Sheet 2 ("bounds in mem or regs graph") raw data edits
# ---------------------------------------------------------------------------

$ws2.Range("K2").Value = "coremk_ch2"

# Row 3 (unsafe, 9 register pairs)
$ws2.Range("K3").Value = 68
$ws2.Range("N3").Value = 72.9
$ws2.Range("O3").Value = 74.6

# Footnote label, same text as Q5 above
$ws2.Range("A48").Value = "UPDATED 20180326"

# ---------------------------------------------------------------------------
# Cosmetic view-state changes (scroll position / selection)
# ---------------------------------------------------------------------------

$ws1.Range("R5").Select()
$ws2.Range("B46").Select()
